$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D), Volume(1h) (E) and Hora (G) columns hold values that look
# numeric/percentage but must stay literal text, matching the source data.
# Pre-format the exact cells we are about to rewrite as Text ("@") so the
# COM layer does not auto-coerce the new literals into numbers/percentages.
# (Ranges are set one contiguous block at a time -- multi-area "A1,B2" refs
# only apply NumberFormat to the first area.)
$ws.Range("D2:D13").NumberFormat = "@"
$ws.Range("D15:D17").NumberFormat = "@"
$ws.Range("D20:D21").NumberFormat = "@"
$ws.Range("D23:D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D39:D42").NumberFormat = "@"
$ws.Range("D44:D50").NumberFormat = "@"
$ws.Range("E2:E18").NumberFormat = "@"
$ws.Range("E20:E21").NumberFormat = "@"
$ws.Range("E23:E27").NumberFormat = "@"
$ws.Range("E39:E50").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "300.21"
$ws.Range("E2").Value = "-0.11%"
$ws.Range("G2").Value = "20"

# Row 3
$ws.Range("D3").Value = "31.81"
$ws.Range("E3").Value = "1.44%"
$ws.Range("G3").Value = "20"

# Row 4
$ws.Range("D4").Value = "5.109"
$ws.Range("E4").Value = "-0.24%"
$ws.Range("G4").Value = "20"

# Row 5
$ws.Range("D5").Value = "0.08219"
$ws.Range("E5").Value = "11.54%"
$ws.Range("G5").Value = "20"

# Row 6
$ws.Range("D6").Value = "2.625"
$ws.Range("E6").Value = "13.97%"
$ws.Range("G6").Value = "20"

# Row 7
$ws.Range("D7").Value = "7.850"
$ws.Range("E7").Value = "-1.36%"
$ws.Range("G7").Value = "20"

# Row 8
$ws.Range("D8").Value = "3.845"
$ws.Range("E8").Value = "1.39%"
$ws.Range("G8").Value = "20"

# Row 9
$ws.Range("D9").Value = "0.9293"
$ws.Range("E9").Value = "1.46%"
$ws.Range("G9").Value = "20"

# Row 10
$ws.Range("D10").Value = "0.1755"
$ws.Range("E10").Value = "2.41%"
$ws.Range("G10").Value = "20"

# Row 11
$ws.Range("D11").Value = "0.07513"
$ws.Range("E11").Value = "-1.73%"
$ws.Range("G11").Value = "20"

# Row 12
$ws.Range("D12").Value = "0.09083"
$ws.Range("E12").Value = "12.04%"
$ws.Range("G12").Value = "20"

# Row 13
$ws.Range("D13").Value = "0.03025"
$ws.Range("E13").Value = "0.11%"
$ws.Range("G13").Value = "20"

# Row 14
$ws.Range("E14").Value = "1.00%"
$ws.Range("G14").Value = "20"

# Row 15
$ws.Range("D15").Value = "0.001508"
$ws.Range("E15").Value = "0.45%"
$ws.Range("G15").Value = "20"

# Row 16
$ws.Range("D16").Value = "0.006050"
$ws.Range("E16").Value = "-2.11%"
$ws.Range("G16").Value = "20"

# Row 17
$ws.Range("D17").Value = "3.612"
$ws.Range("E17").Value = "4.25%"
$ws.Range("G17").Value = "20"

# Row 18
$ws.Range("E18").Value = "2.81%"
$ws.Range("G18").Value = "20"

# Row 19
$ws.Range("G19").Value = "20"

# Row 20
$ws.Range("D20").Value = "0.1346"
$ws.Range("E20").Value = "0.73%"
$ws.Range("G20").Value = "20"

# Row 21
$ws.Range("D21").Value = "3.916"
$ws.Range("E21").Value = "-15.78%"
$ws.Range("G21").Value = "20"

# Row 22
$ws.Range("G22").Value = "20"

# Row 23
$ws.Range("D23").Value = "0.04614"
$ws.Range("E23").Value = "-0.74%"
$ws.Range("G23").Value = "20"

# Row 24
$ws.Range("D24").Value = "0.001245"
$ws.Range("E24").Value = "1.62%"
$ws.Range("G24").Value = "20"

# Row 25
$ws.Range("D25").Value = "0.004556"
$ws.Range("E25").Value = "1.60%"
$ws.Range("G25").Value = "20"

# Row 26
$ws.Range("E26").Value = "-7.74%"
$ws.Range("G26").Value = "20"

# Row 27
$ws.Range("D27").Value = "0.0003401"
$ws.Range("E27").Value = "81.71%"
$ws.Range("G27").Value = "20"

# Row 28
$ws.Range("G28").Value = "20"

# Row 29
$ws.Range("G29").Value = "20"

# Row 30
$ws.Range("G30").Value = "20"

# Row 31
$ws.Range("G31").Value = "20"

# Row 32
$ws.Range("G32").Value = "20"

# Row 33
$ws.Range("G33").Value = "20"

# Row 34
$ws.Range("G34").Value = "20"

# Row 35
$ws.Range("G35").Value = "20"

# Row 36
$ws.Range("G36").Value = "20"

# Row 37
$ws.Range("G37").Value = "20"

# Row 38
$ws.Range("G38").Value = "20"

# Row 39
$ws.Range("D39").Value = "0.01779"
$ws.Range("E39").Value = "2.81%"
$ws.Range("G39").Value = "20"

# Row 40
$ws.Range("D40").Value = "0.04613"
$ws.Range("E40").Value = "2.64%"
$ws.Range("G40").Value = "20"

# Row 41
$ws.Range("D41").Value = "0.006880"
$ws.Range("E41").Value = "-4.82%"
$ws.Range("G41").Value = "20"

# Row 42
$ws.Range("D42").Value = "0.1381"
$ws.Range("E42").Value = "2.52%"
$ws.Range("G42").Value = "20"

# Row 43
$ws.Range("E43").Value = "0.41%"
$ws.Range("G43").Value = "20"

# Row 44
$ws.Range("D44").Value = "0.009757"
$ws.Range("E44").Value = "-8.96%"
$ws.Range("G44").Value = "20"

# Row 45
$ws.Range("D45").Value = "0.00006161"
$ws.Range("E45").Value = "-1.94%"
$ws.Range("G45").Value = "20"

# Row 46
$ws.Range("D46").Value = "0.00000000748"
$ws.Range("E46").Value = "-0.28%"
$ws.Range("G46").Value = "20"

# Row 47
$ws.Range("B47").Value = "CoinbaseStockToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D47").Value = "0.008377"
$ws.Range("E47").Value = "-16.25%"
$ws.Range("G47").Value = "20"

# Row 48
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "0.7900"
$ws.Range("E48").Value = "-2.28%"
$ws.Range("G48").Value = "20"

# Row 49
$ws.Range("D49").Value = "0.00002094"
$ws.Range("E49").Value = "-0.28%"
$ws.Range("G49").Value = "20"

# Row 50
$ws.Range("D50").Value = "0.0001995"
$ws.Range("E50").Value = "-0.21%"
$ws.Range("G50").Value = "20"

# Row 51
$ws.Range("G51").Value = "20"
